$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("general")
$ws.Range("B3").Value = 181.9610138849955
$ws.Range("B4").Value = 0.0149998664855957
$ws.Range("B6").Value = 32.48101388499555
$ws.Range("B9").Value = 149.48
$ws.Range("B10").Value = 0

$ws = $wb.Worksheets.Item("x")
$ws.Range("B3").Value = 5
$ws.Range("B7").Value = 13
$ws.Range("B9").Value = 6
$ws.Range("B13").Value = 11

$ws = $wb.Worksheets.Item("U")
$ws.Range("B11").Value = 3

$ws = $wb.Worksheets.Item("TBar")
$ws.Range("B3").Value = 22.61192465059683
$ws.Range("B5").Value = 24.16886835983306
$ws.Range("B7").Value = 20
$ws.Range("B8").Value = 23.82001308898281
$ws.Range("B9").Value = 22.01159140980468
$ws.Range("B12").Value = 20
$ws.Range("B13").Value = 27.87444125446785
$ws.Range("B15").Value = 26.48759645946009

$ws = $wb.Worksheets.Item("Q")
$ws.Range("C7").Value = 250.9699999999993
$ws.Range("C8").Value = 260.9899999999993
$ws.Range("C9").Value = 252.9749999999993
$ws.Range("C10").Value = 269.58
$ws.Range("C11").Value = 250.5749999999993
$ws.Range("C12").Value = 332.445
$ws.Range("C13").Value = 341.67
$ws.Range("C14").Value = 337.69
$ws.Range("C15").Value = 350.0150000000007
$ws.Range("C16").Value = 337.92
$ws.Range("C17").Value = 224.1799999999995
$ws.Range("C18").Value = 224.6649999999995
$ws.Range("C19").Value = 201.1149999999995
$ws.Range("C20").Value = 218.9699999999995
$ws.Range("C21").Value = 207.1049999999995
$ws.Range("C27").Value = 69.88000000000022
$ws.Range("C28").Value = 76.32000000000022
$ws.Range("C29").Value = 72.41000000000021
$ws.Range("C30").Value = 72.64500000000022
$ws.Range("C31").Value = 72.89000000000021
$ws.Range("C32").Value = 177.26
$ws.Range("C33").Value = 188.5600000000001
$ws.Range("C34").Value = 166.2
$ws.Range("C35").Value = 181.48
$ws.Range("C36").Value = 167.6700000000001
$ws.Range("C37").Value = 141.0250000000001
$ws.Range("C38").Value = 143.4
$ws.Range("C39").Value = 139.7050000000001
$ws.Range("C40").Value = 150.4249999999992
$ws.Range("C41").Value = 134.7700000000002
$ws.Range("C43").Value = 159.2149999999989
$ws.Range("C46").Value = 139.7449999999989
$ws.Range("C51").Value = 224.4749999999994
$ws.Range("C52").Value = 57.95
$ws.Range("C53").Value = 58.67999999999927
$ws.Range("C54").Value = 61.72999999999927
$ws.Range("C55").Value = 60.65499999999928
$ws.Range("C56").Value = 52.91499999999927
$ws.Range("C57").Value = 226.0399999999994
$ws.Range("C58").Value = 247.1799999999994
$ws.Range("C59").Value = 221.8549999999994
$ws.Range("C60").Value = 238.4549999999994
$ws.Range("C61").Value = 224.4749999999994
$ws.Range("C62").Value = 332.445
$ws.Range("C63").Value = 341.67
$ws.Range("C64").Value = 337.69
$ws.Range("C65").Value = 350.0150000000007
$ws.Range("C66").Value = 337.92
$ws.Range("C67").Value = 224.1799999999995
$ws.Range("C68").Value = 224.6649999999995
$ws.Range("C69").Value = 201.1149999999995
$ws.Range("C70").Value = 218.9699999999995
$ws.Range("C71").Value = 207.1049999999995

$ws = $wb.Worksheets.Item("R")
$ws.Range("C7").Value = 7.445
$ws.Range("C8").Value = 16.67
$ws.Range("C9").Value = 12.69
$ws.Range("C10").Value = 25.015
$ws.Range("C11").Value = 12.92

$ws = $wb.Worksheets.Item("L")
$ws.Range("C47").Value = 0
$ws.Range("C48").Value = 0
$ws.Range("C49").Value = 0
$ws.Range("C50").Value = 0
$ws.Range("C51").Value = 0
